$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "10.00", "2.060") keep their exact original formatting instead
# of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.173.33"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "1.905.00"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "327.29"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("D8").Value = "0.3954"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "46.81"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "0.07964"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").Value = "22.33"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "1.917.63"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "7.143"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "5.793"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "0.06958"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "88.81"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D20").Value = "17.19"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "29.203.08"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "5.355"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "2.143.55"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "2.060"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "156.67"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").Value = "19.57"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "5.877"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "1.999"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "119.64"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "0.09441"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "0.9221"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "5.357"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "3.258"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Value = "0.05848"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "1.172"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02104"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.979"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").Value = "0.5752"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "0.1813"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "10.00"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "12.10"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("D45").Value = "0.5427"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").Value = "2.212"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").Value = "0.07092"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "1.887"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("D49").Value = "2.559"
$ws.Range("E49").Value = "  +6.21%  "
$ws.Range("D50").Value = "112.11"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "1.073"
$ws.Range("E51").Value = "  -5.53%  "

# Restore default style on column D so no stray number-format style
# reference is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"

